$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D-column cells that look like plain numbers need a leading apostrophe
# so Excel keeps them as text (matching the original inlineStr cells).

$ws.Range("D2").Value = '25.959.71'
$ws.Range("E2").Value = '  +0.47%  '
$ws.Range("D3").Value = '1.737.58'
$ws.Range("E3").Value = '  +0.04%  '
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = "'246.99"
$ws.Range("E5").Value = '  +4.24%  '
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = '  -0.04%  '
$ws.Range("D7").Value = "'0.5038"
$ws.Range("E7").Value = '  -1.40%  '
$ws.Range("D8").Value = "'0.2726"
$ws.Range("E8").Value = '  -0.22%  '
$ws.Range("D9").Value = "'0.06182"
$ws.Range("E9").Value = '  +1.27%  '
$ws.Range("E10").Value = '  +1.30%  '
$ws.Range("D11").Value = '1.737.44'
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("D12").Value = "'0.6569"
$ws.Range("E12").Value = '  +3.18%  '
$ws.Range("D13").Value = "'15.26"
$ws.Range("E13").Value = '  +2.59%  '
$ws.Range("D14").Value = "'4.756"
$ws.Range("E14").Value = '  +3.75%  '
$ws.Range("D15").Value = "'77.60"
$ws.Range("E15").Value = '  +0.70%  '
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = '  -0.07%  '
$ws.Range("D17").Value = "'0.9998"
$ws.Range("E17").Value = '  -0.08%  '
$ws.Range("D18").Value = '25.981.08'
$ws.Range("E18").Value = '  +0.53%  '
$ws.Range("D19").Value = "'11.88"
$ws.Range("E19").Value = '  +1.58%  '
$ws.Range("D20").Value = "'0.000006823"
$ws.Range("E20").Value = '  +1.39%  '
$ws.Range("E21").Value = '  +8.63%  '
$ws.Range("D22").Value = '1.962.20'
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = "'8.806"
$ws.Range("E23").Value = '  +1.67%  '
$ws.Range("D24").Value = "'5.451"
$ws.Range("E24").Value = '  +4.45%  '
$ws.Range("D25").Value = "'134.64"
$ws.Range("E25").Value = '  -3.09%  '
$ws.Range("D26").Value = "'15.26"
$ws.Range("E26").Value = '  +1.07%  '
$ws.Range("D27").Value = "'1.456"
$ws.Range("E27").Value = '  -3.81%  '
$ws.Range("D28").Value = "'1.791"
$ws.Range("E28").Value = '  +2.21%  '
$ws.Range("D29").Value = "'105.43"
$ws.Range("E29").Value = '  -0.06%  '
$ws.Range("D30").Value = "'3.994"
$ws.Range("E30").Value = '  +0.06%  '
$ws.Range("D31").Value = "'0.08136"
$ws.Range("E31").Value = '  -2.08%  '
$ws.Range("D32").Value = "'3.722"
$ws.Range("E32").Value = '  +1.96%  '
$ws.Range("D33").Value = "'0.04767"
$ws.Range("E33").Value = '  +4.59%  '
$ws.Range("D34").Value = "'2.654"
$ws.Range("E34").Value = '  -0.43%  '
$ws.Range("D35").Value = "'0.9976"
$ws.Range("E35").Value = '  +1.36%  '
$ws.Range("D36").Value = "'0.6119"
$ws.Range("E36").Value = '  -0.59%  '
$ws.Range("D37").Value = "'2.738"
$ws.Range("E37").Value = '  +1.98%  '
$ws.Range("D38").Value = "'0.01613"
$ws.Range("E38").Value = '  +0.13%  '
$ws.Range("D39").Value = "'0.8668"
$ws.Range("E39").Value = '  +18.00%  '
$ws.Range("D40").Value = "'1.956"
$ws.Range("E40").Value = '  +1.74%  '
$ws.Range("D41").Value = "'0.9996"
$ws.Range("E41").Value = '  -0.06%  '
$ws.Range("D42").Value = "'100.86"
$ws.Range("E42").Value = '  +3.11%  '
$ws.Range("D43").Value = "'0.3940"
$ws.Range("E43").Value = '  +2.87%  '
$ws.Range("D44").Value = "'5.029"
$ws.Range("E44").Value = '  +1.72%  '
$ws.Range("E45").Value = '  +5.58%  '
$ws.Range("D46").Value = "'6.367"
$ws.Range("E46").Value = '  +3.51%  '
$ws.Range("D47").Value = "'55.73"
$ws.Range("E47").Value = '  +1.74%  '
$ws.Range("D48").Value = "'0.05280"
$ws.Range("E48").Value = '  +0.39%  '
$ws.Range("D49").Value = "'30.95"
$ws.Range("E49").Value = '  +1.70%  '
$ws.Range("D50").Value = "'0.3498"
$ws.Range("E50").Value = '  +2.80%  '
$ws.Range("D51").Value = "'7.650"
$ws.Range("E51").Value = '  +0.74%  '
